# Correcting the npv saving. The loans are changed for the restpayment.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate()

# start_year_dismantling: 4 -> 10
$ws.Range("B12").Value = 10

# maximum_investment_capacity_per_year: 600 -> 150
$ws.Range("B13").Value = 150

# typeofProfitforPastHorizon: totalProfits -> none
$ws.Range("B14").Value = "none"

# Widen column C to fit the longer parameter descriptions (was 54.09)
$ws.Columns.Item(3).ColumnWidth = 69.27

# Update the selected cell to reflect where the edits were made
$ws.Range("C8").Select()
